$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CREATED_BY (column D) for both data rows from "05582" to the email address.
$ws.Range("D2").Value = "SWAPNOTORY49@GMAIL.COM"
$ws.Range("D3").Value = "SWAPNOTORY49@GMAIL.COM"

# Update CREATED_AT (column E) timestamps for both rows.
$ws.Range("E2").Value = 45312.7411574074
$ws.Range("E3").Value = 45312.7411574074

# Ensure column D is wide enough to fit the new, longer value (target stored width ~31.36).
$ws.Columns.Item(4).ColumnWidth = 30.5
